$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K column (G) values per row, replacing old "Strike#" derived counts
# with the regenerated K values (calc + write s_vals).
$kValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 1
    8 = 2
    9 = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    21 = 0
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 2
    47 = 0
    48 = 2
    49 = 0
    50 = 1
    51 = 0
    52 = 0
    53 = 1
    54 = 0
    55 = 1
    56 = 0
    58 = 1
    59 = 0
    60 = 1
    61 = 1
    62 = 0
    63 = 3
    64 = 1
    65 = 1
    66 = 2
    67 = 2
    68 = 0
    69 = 1
    70 = 2
    71 = 2
    72 = 1
    73 = 2
    74 = 0
    75 = 1
    76 = 0
    77 = 0
    78 = 1
    80 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
